$d = $word.ActiveDocument
$wNS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ------------------------------------------------------------------
# 1) "Permittee: {d.permitee}" paragraph
#    - drop the both-justify alignment
#    - fix the field name typo d.permitee -> d.permittee
# ------------------------------------------------------------------
$permitteeFound = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Permittee: *") {
        $p.Format.Alignment = 0  # wdAlignParagraphLeft -> removes <w:jc w:val="both"/>
        $p.Range.Find.Execute("d.permitee", $true, $false, $false, $false, $false, $true, 1, $false, "d.permittee", 2)
        $permitteeFound = $true
        break
    }
}
Write-Host "permittee paragraph fixed: $permitteeFound"

# ------------------------------------------------------------------
# 2) "Address of Permittee: {d.permitee_mailing_address}" paragraph
#    -> replaced by a 2 column, border-less table.
# ------------------------------------------------------------------

# Make sure the "TableGrid" table style exists in the document with the
# right priority/spacing so the table we insert below can reference it.
$seed = $d.Tables.Add($d.Range($d.Content.End - 1, $d.Content.End - 1), 1, 1)
$seed.Style = "Table Grid"
$tgStyle = $d.Styles("Table Grid")
$tgStyle.Priority = 39
$tgStyle.ParagraphFormat.SpaceAfter = 0
$tgStyle.ParagraphFormat.LineSpacingRule = 0
$seed.Delete()

$addrParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Address of Permittee:*") {
        $addrParagraph = $p
        break
    }
}

if ($addrParagraph -ne $null) {
    $insAt = $d.Range($addrParagraph.Range.Start, $addrParagraph.Range.Start)

    $rPrPlain = '<w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'
    $rPrBold = '<w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'
    $cellPPr = '<w:pPr><w:spacing w:line="300" w:lineRule="exact"/><w:rPr><w:b/><w:i/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>'
    $tcBorders = '<w:tcBorders><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/></w:tcBorders>'

    $cell1 = '<w:tc><w:tcPr><w:tcW w:w="2029" w:type="dxa"/>' + $tcBorders + '</w:tcPr><w:p>' + $cellPPr + '<w:r>' + $rPrPlain + '<w:t>Address of Permittee:</w:t></w:r></w:p></w:tc>'

    $cell2Runs = @(
        ('<w:r>' + $rPrBold + '<w:t>{</w:t></w:r>'),
        '<w:proofErr w:type="spellStart"/>',
        ('<w:r>' + $rPrBold + '<w:t>d.permit</w:t></w:r>'),
        ('<w:r>' + $rPrBold + '<w:t>t</w:t></w:r>'),
        ('<w:r>' + $rPrBold + '<w:t>ee_mailing_address</w:t></w:r>'),
        ('<w:r>' + $rPrBold + '<w:t>:convCRLF</w:t></w:r>'),
        '<w:proofErr w:type="spellEnd"/>',
        ('<w:r>' + $rPrBold + '<w:t>()</w:t></w:r>'),
        ('<w:r>' + $rPrBold + '<w:t>}</w:t></w:r>')
    ) -join ''

    $cell2 = '<w:tc><w:tcPr><w:tcW w:w="4274" w:type="dxa"/>' + $tcBorders + '</w:tcPr><w:p>' + $cellPPr + $cell2Runs + '</w:p></w:tc>'

    $tblPr = '<w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblCellMar><w:left w:w="0" w:type="dxa"/><w:right w:w="0" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr>'
    $tblGrid = '<w:tblGrid><w:gridCol w:w="2029"/><w:gridCol w:w="4274"/></w:tblGrid>'

    $tableXml = '<w:tbl xmlns:w="' + $wNS + '">' + $tblPr + $tblGrid + '<w:tr>' + $cell1 + $cell2 + '</w:tr></w:tbl>'

    $insAt.InsertXML($tableXml)
    Write-Host "table inserted"

    # Remove the original (now orphaned) "Address of Permittee" paragraph.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "Address of Permittee: {d.permitee_mailing_address}*") {
            $p.Range.Delete()
            break
        }
    }
    Write-Host "orphan paragraph removed"
}

# ------------------------------------------------------------------
# 3) Merge split runs in the "The amended Mines Act..." paragraph.
# ------------------------------------------------------------------
$bigParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The amended Mines Act Permit*") {
        $bigParagraph = $p
        break
    }
}

if ($bigParagraph -ne $null) {
    $rng1 = $bigParagraph.Range
    $rng1.Find.Execute("The amended ", $true, $false, $false, $false, $false, $true, 1, $false, "The amended ", 2)

    $rng2 = $bigParagraph.Range
    $mid = " contain the requirements of the Ministry of Energy, Mines and Petroleum Resources for reclamation, including reclamation securities. "
    $rng2.Find.Execute($mid, $true, $false, $false, $false, $false, $true, 1, $false, $mid, 2)

    Write-Host "merged runs in amended paragraph"
}
